$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timestamps = @(
    "2021-10-05 10:50:50.176419",
    "2021-10-05 10:50:50.176432",
    "2021-10-05 10:50:50.176436",
    "2021-10-05 10:50:50.176439",
    "2021-10-05 10:50:50.176442",
    "2021-10-05 10:50:50.176446",
    "2021-10-05 10:50:50.176449",
    "2021-10-05 10:50:50.176452",
    "2021-10-05 10:50:50.176455",
    "2021-10-05 10:50:50.176459",
    "2021-10-05 10:50:50.176462",
    "2021-10-05 10:50:50.176465",
    "2021-10-05 10:50:50.176468",
    "2021-10-05 10:50:50.176471",
    "2021-10-05 10:50:50.176474",
    "2021-10-05 10:50:50.176477",
    "2021-10-05 10:50:50.176481",
    "2021-10-05 10:50:50.176484",
    "2021-10-05 10:50:50.176487",
    "2021-10-05 10:50:50.176490",
    "2021-10-05 10:50:50.176493",
    "2021-10-05 10:50:50.176496",
    "2021-10-05 10:50:50.176499",
    "2021-10-05 10:50:50.176502",
    "2021-10-05 10:50:50.176506",
    "2021-10-05 10:50:50.176509",
    "2021-10-05 10:50:50.176512",
    "2021-10-05 10:50:50.176515",
    "2021-10-05 10:50:50.176518",
    "2021-10-05 10:50:50.176521",
    "2021-10-05 10:50:50.176524",
    "2021-10-05 10:50:50.176528",
    "2021-10-05 10:50:50.176531",
    "2021-10-05 10:50:50.176534",
    "2021-10-05 10:50:50.176537",
    "2021-10-05 10:50:50.176540",
    "2021-10-05 10:50:50.176543",
    "2021-10-05 10:50:50.176547",
    "2021-10-05 10:50:50.176550",
    "2021-10-05 10:50:50.176553",
    "2021-10-05 10:50:50.176557",
    "2021-10-05 10:50:50.176560",
    "2021-10-05 10:50:50.176563",
    "2021-10-05 10:50:50.176566",
    "2021-10-05 10:50:50.176569",
    "2021-10-05 10:50:50.176572",
    "2021-10-05 10:50:50.176576",
    "2021-10-05 10:50:50.176579",
    "2021-10-05 10:50:50.176582"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
